$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 28
$ws.Columns.Item(2).ColumnWidth = 22
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 22
$ws.Columns.Item(6).ColumnWidth = 28

# Header row
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$asesor = "RIOS CARRION ANGEL BENIGNO"

$ws.Range("A2").Value = $asesor
$ws.Range("B2").Value = "240X120 PORCELANATO"
$ws.Range("C2").Value = 1041.16
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1041.16
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = $asesor
$ws.Range("B3").Value = "240X80 PORCELANATO"
$ws.Range("C3").Value = 8668.91
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 8668.91
$ws.Range("F3").Value = 0

$ws.Range("A4").Value = $asesor
$ws.Range("B4").Value = "FREGADEROS DE COCINA"
$ws.Range("C4").Value = 372.993863046034
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 372.993863046034
$ws.Range("F4").Value = 0

$ws.Range("A5").Value = $asesor
$ws.Range("B5").Value = "GRANITO"
$ws.Range("C5").Value = 238.32
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 238.32
$ws.Range("F5").Value = 0

$ws.Range("A6").Value = $asesor
$ws.Range("B6").Value = "GRIFERIAS"
$ws.Range("C6").Value = 106.82
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 106.82
$ws.Range("F6").Value = 0

$ws.Range("A7").Value = $asesor
$ws.Range("B7").Value = "INODOROS"
$ws.Range("C7").Value = 800
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 800
$ws.Range("F7").Value = 0

$ws.Range("A8").Value = $asesor
$ws.Range("B8").Value = "LAVABOS"
$ws.Range("C8").Value = 625
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 625
$ws.Range("F8").Value = 0

$ws.Range("A9").Value = $asesor
$ws.Range("B9").Value = "LED"
$ws.Range("C9").Value = 300
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 300
$ws.Range("F9").Value = 0

$ws.Range("A10").Value = $asesor
$ws.Range("B10").Value = "NO RESURTIBLES"
$ws.Range("C10").Value = 650.25
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 650.25
$ws.Range("F10").Value = 0

$ws.Range("A11").Value = $asesor
$ws.Range("B11").Value = "OTROS"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

$ws.Range("A12").Value = $asesor
$ws.Range("B12").Value = "PANELES DECORATIVOS"
$ws.Range("C12").Value = 350
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 350
$ws.Range("F12").Value = 0

$ws.Range("A13").Value = $asesor
$ws.Range("B13").Value = "PANELES PU"
$ws.Range("C13").Value = 230
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 230
$ws.Range("F13").Value = 0

$ws.Range("A14").Value = $asesor
$ws.Range("B14").Value = "PANELES PVC"
$ws.Range("C14").Value = 483
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 483
$ws.Range("F14").Value = 0

$ws.Range("A15").Value = $asesor
$ws.Range("B15").Value = "PIEDRA SINTERIZADA"
$ws.Range("C15").Value = 2501.01
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 2501.01
$ws.Range("F15").Value = 0

$ws.Range("A16").Value = $asesor
$ws.Range("B16").Value = "PORCELANATO"
$ws.Range("C16").Value = 28209.84
$ws.Range("D16").Value = -22.68
$ws.Range("E16").Value = 28232.52
$ws.Range("F16").Value = -0.0008039747832671153

$ws.Range("A17").Value = $asesor
$ws.Range("B17").Value = "PUERTAS DE SEGURIDAD"
$ws.Range("C17").Value = 342
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 342
$ws.Range("F17").Value = 0

$ws.Range("A18").Value = $asesor
$ws.Range("B18").Value = "SAL SOLUBLE"
$ws.Range("C18").Value = 2300
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 2300
$ws.Range("F18").Value = 0

$ws.Range("B19").Value = "TOTAL"
$ws.Range("C19").Value = 47219.30386304604
$ws.Range("D19").Value = -22.68
$ws.Range("E19").Value = 47241.98386304604
$ws.Range("F19").Value = -0.0004803120364878872

# Number formats and alignment
$currencyFmt = """$""#,##0.00"
$ws.Range("C2:E19").NumberFormat = $currencyFmt
$ws.Range("F2:F19").NumberFormat = "0.00%"
$ws.Range("B19").HorizontalAlignment = -4152

# Move the new sheet to the end (after "VENTA MENSUAL"), must be last
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)
